$d = $word.ActiveDocument

# 1. "Engineered Docker setup for local dev of custom company theme"
#    -> full rewrite describing contractor onboarding impact
$d.Content.Find.Execute(
    "Engineered Docker setup for local dev of custom company theme",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Engineered Docker setup to onboard contractors in a single afternoon, replacing multi-day manual steps with automated scripts and Docker Compose",
    2) | Out-Null

# 2. "Facilitated company website deployments..." -> "Deployed company websites..."
$d.Content.Find.Execute(
    "Facilitated company website deployments via SSH and SpinupWP, including server and Cloudflare cache handling",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployed company websites via SSH and SpinupWP, managing servers and Cloudflare cache",
    2) | Out-Null

# 3. Swap word order: "manual 10-minute" -> "10-minute manual"
#    (use a narrow Range so the untouched surrounding runs are left alone)
$rngManual = $d.Content
$rngManual.Find.Execute("manual 10-minute") | Out-Null
if ($rngManual.Find.Found) {
    $rngManual.Text = "10-minute manual"
}

# 4. Insert "securely " before "run accessibility evaluations"
$rngSecure = $d.Content
$rngSecure.Find.Execute("run accessibility evaluations") | Out-Null
if ($rngSecure.Find.Found) {
    $rngSecure.Text = "securely run accessibility evaluations"
}

# 6. Flask game bullet: "for more than" -> "to engage over", add comma in "tight, two-week"
$d.Content.Find.Execute(
    "Deployed Flask-based multiple-choice game for more than 500 users at conference and convention events, built and launched under tight two-week deadline",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Deployed Flask-based multiple-choice game to engage over 500 users at conference and convention events, built and launched under tight, two-week deadline",
    2) | Out-Null

# 7. Languages line: HTML/CSS/JS -> HTML, CSS, JavaScript
$d.Content.Find.Execute(
    ", Java, HTML/CSS/JS",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", Java, HTML, CSS, JavaScript",
    2) | Out-Null

# 8. Frameworks: Django, Flask -> Frameworks: Flask, Django
$d.Content.Find.Execute(
    "Frameworks: Django, Flask",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Frameworks: Flask, Django",
    2) | Out-Null

# 11. "Hosting & CMS: " -> "CMS & Hosting: "
$d.Content.Find.Execute(
    "Hosting & CMS: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CMS & Hosting: ",
    2) | Out-Null
